# Change container length/width values (swap columns B and C for rows 2-6)
# and update the active selection, as part of passing container info to visualizer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 235
$ws.Range("C2").Value = 590

$ws.Range("B3").Value = 235
$ws.Range("C3").Value = 1204

$ws.Range("B4").Value = 235
$ws.Range("C4").Value = 1204

$ws.Range("B5").Value = 244
$ws.Range("C5").Value = 1359

$ws.Range("B6").Value = 248
$ws.Range("C6").Value = 1360

# Update selection to match new active cell
$ws.Range("C11").Select()
